# Auto-generated edit script applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names / links) ---
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"

# --- Numeric-looking text cells (price / volume%) kept as Text ---
# Force text number format so values like "242.61" or "-1.10%" are stored as literal strings,
# matching the original inlineStr text cells rather than being parsed into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.10%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "9.94%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.102"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.60%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05651"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.13%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.499"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.37%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8256"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.07%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8686"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.31%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.01005"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1,578.13%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1331"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.12%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06916"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.19%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02857"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.02%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09375"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.15%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001516"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.59%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04167"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-9.27%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006065"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.57%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.521"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.19%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.025"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.39%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.219"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.65%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03244"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.80%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.10%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.613"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.33%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.00%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001211"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.58%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004442"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.86%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001179"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "22.91%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001403"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "0.57%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03706"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.75%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005918"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.73%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.17%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002311"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.69%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009555"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "15.55%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005096"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.50%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.00%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1050"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-3.67%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002436"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-4.38%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.00%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
